# "updated Week 1 lession plan"
#
# Applies the text edits described by the diff. Every change below is a
# Find/Replace against a phrase that is unique within the document, so
# each call only ever touches the one paragraph it is meant for. The
# replacement text keeps the same run formatting as the text it
# replaces (Word COM's Find/Replace does not alter character
# formatting), which matches every hunk in the diff -- each hunk only
# splits a run into pieces that all carry identical rPr, i.e. no visual
# formatting actually changes anywhere.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute(
        $old,    # FindText
        $true,   # MatchCase
        $false,  # MatchWholeWord
        $false,  # MatchWildcards
        $false,  # MatchSoundsLike
        $false,  # MatchAllWordForms
        $true,   # Forward
        1,       # Wrap (wdFindContinue)
        $false,  # Format
        $new,    # ReplaceWith
        2        # Replace (wdReplaceAll)
    ) | Out-Null
}

# 1. Title gains "R " before "Markdown File".
Replace-Text "Set Up Your R Project & Markdown File" "Set Up Your R Project & R Markdown File"

# 2. "Just open the RStudio..." -> "Open the RStudio..."
Replace-Text "Just open the RStudio application on your computer." "Open the RStudio application on your computer."

# 3. Drop "some " before "panes".
Replace-Text "You should see an empty console window and some panes for files, plots, etc." "You should see an empty console window and panes for files, plots, etc."

# 4. Replace the em dash aside with the "R course folder" phrasing.
$emDash = [char]0x2014
$lq = [char]0x201C
$rq = [char]0x201D
Replace-Text (" and pick a folder" + $emDash + "like " + $lq + "Documents" + $rq + " or your desktop).") `
             " and pick a folder inside your R course folder that you already set up)."

# 5. "Week 1 Practice" -> "Week 1 Skills Learning"
Replace-Text "Week 1 Practice" "Week 1 Skills Learning"

# 6. Shorten the Author hint.
Replace-Text " (you can leave blank or add your name)" " (add your name)"

# 7. Rename the sample filename.
Replace-Text "week1_practice.Rmd" "week1_skills_learning.Rmd"

# 8. "in the white text parts." -> "in between code chunks."
Replace-Text " in the white text parts." " in between code chunks."
